$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("F3").Value = 800
$ws.Range("F4").Value = 500
$ws.Range("F8").Value = 800
$ws.Range("F9").Value = 450
$ws.Range("F10").Value = 800
$ws.Range("F11").Value = 800
$ws.Range("F12").Value = 800
$ws.Range("F13").Value = 800
$ws.Range("F15").Value = 800
$ws.Range("F16").Value = 800
$ws.Range("F17").Value = 800
$ws.Range("F18").Value = 800
$ws.Range("F19").Value = 800
$ws.Range("F20").Value = 800
$ws.Range("F21").Value = 500
$ws.Range("F22").Value = 800
